# Refactor currency conversion sheet: split "foreign_amount" into explicit
# "source_amount" / "source_fees" / "source_currency" and
# "target_amount" / "target_fees" / "target_currency" columns, and make
# "currency_conversions" the active sheet.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "currency_conversions": rebuild header row and data row with the
# new column layout:
#   A date | B source_amount | C source_fees | D source_currency |
#   E target_amount | F target_fees | G target_currency | H comment
# ---------------------------------------------------------------------
$ws4 = $wb.Worksheets.Item("currency_conversions")

# header text (write B1 last so new shared strings are appended in the
# same order as target_amount, target_fees, source_amount)
$ws4.Range("A1").Value2 = "date"
$ws4.Range("C1").Value2 = "source_fees"
$ws4.Range("D1").Value2 = "source_currency"
$ws4.Range("E1").Value2 = "target_amount"
$ws4.Range("F1").Value2 = "target_fees"
$ws4.Range("G1").Value2 = "target_currency"
$ws4.Range("H1").Value2 = "comment"
$ws4.Range("B1").Value2 = "source_amount"

# give the newly-extended header cells (G1, H1) the same look as the
# other header cells (bold, centered)
$ws4.Range("A1").Copy()
$ws4.Range("G1").PasteSpecial(-4122)
$ws4.Range("H1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# B1 ("source_amount") uses a plain bold style (no center/border)
$ws4.Range("B1").ClearFormats()
$ws4.Range("B1").Font.Name = "Calibri"
$ws4.Range("B1").Font.FontStyle = "Bold"
$ws4.Range("B1").Font.Size = 11
$ws4.Range("B1").Font.ThemeColor = 1

# data row
$ws4.Range("B2").Value2 = -1
$ws4.Range("E2").Value2 = 150
$ws4.Range("F2").Value2 = 0
$ws4.Range("G2").Value2 = "USD"

# ---------------------------------------------------------------------
# Sheet "rsu": header row re-styled (no effective content change)
# ---------------------------------------------------------------------
$ws6 = $wb.Worksheets.Item("rsu")
$ws6.Range("A1").Copy()
$ws6.Range("A1:H1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# ---------------------------------------------------------------------
# Sheet "espp": header row re-styled (no effective content change)
# ---------------------------------------------------------------------
$ws7 = $wb.Worksheets.Item("espp")
$ws4.Range("A1").Copy()
$ws7.Range("A1:G1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# ---------------------------------------------------------------------
# Make "currency_conversions" the active / selected sheet
# ---------------------------------------------------------------------
$ws4.Activate()
